# Femacal de La Calera - Damasco: weekly update
# Inserts a new "Patterson" price block (week of 2023-01-06) right after the
# existing "Patterson" rows at 170-171, pushing every later block down by
# three rows. Because the insert happens above them, the original last block
# (rows 213-215, Castle Brite @ 2022-12-16) automatically lands on rows
# 216-218 with its data/formatting intact - no further action is required
# for the tail of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above row 172 (formatting of row 172, incl. the date
# style on column D, is carried down into the new rows by Excel).
$ws.Rows.Item(172).Resize(3).Insert()

# Row 172: Patterson / Especial
$ws.Range("A172").Value = 3
$ws.Range("B172").Value = "Femacal de La Calera"
$ws.Range("C172").Value = "Coquimbo"
$ws.Range("D172").Value = 44932
$ws.Range("E172").Value = 5
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100103
$ws.Range("H172").Value = "Frutos de hueso (carozo)"
$ws.Range("I172").Value = 100103003
$ws.Range("J172").Value = "Damasco"
$ws.Range("K172").Value = "Patterson"
$ws.Range("L172").Value = "Especial"
$ws.Range("M172").Value = 56
$ws.Range("N172").Value = 14000
$ws.Range("O172").Value = 14000
$ws.Range("P172").Value = 14000
$ws.Range("Q172").Value = "$/caja 15 kilos"
$ws.Range("R172").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S172").Value = 933
$ws.Range("T172").Value = 15

# Row 173: Patterson / Primera
$ws.Range("A173").Value = 3
$ws.Range("B173").Value = "Femacal de La Calera"
$ws.Range("C173").Value = "Coquimbo"
$ws.Range("D173").Value = 44932
$ws.Range("E173").Value = 5
$ws.Range("F173").Value = "Fruta"
$ws.Range("G173").Value = 100103
$ws.Range("H173").Value = "Frutos de hueso (carozo)"
$ws.Range("I173").Value = 100103003
$ws.Range("J173").Value = "Damasco"
$ws.Range("K173").Value = "Patterson"
$ws.Range("L173").Value = "Primera"
$ws.Range("M173").Value = 60
$ws.Range("N173").Value = 12000
$ws.Range("O173").Value = 12000
$ws.Range("P173").Value = 12000
$ws.Range("Q173").Value = "$/caja 15 kilos"
$ws.Range("R173").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S173").Value = 800
$ws.Range("T173").Value = 15

# Row 174: Patterson / Segunda
$ws.Range("A174").Value = 3
$ws.Range("B174").Value = "Femacal de La Calera"
$ws.Range("C174").Value = "Coquimbo"
$ws.Range("D174").Value = 44932
$ws.Range("E174").Value = 5
$ws.Range("F174").Value = "Fruta"
$ws.Range("G174").Value = 100103
$ws.Range("H174").Value = "Frutos de hueso (carozo)"
$ws.Range("I174").Value = 100103003
$ws.Range("J174").Value = "Damasco"
$ws.Range("K174").Value = "Patterson"
$ws.Range("L174").Value = "Segunda"
$ws.Range("M174").Value = 60
$ws.Range("N174").Value = 10000
$ws.Range("O174").Value = 10000
$ws.Range("P174").Value = 10000
$ws.Range("Q174").Value = "$/caja 15 kilos"
$ws.Range("R174").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S174").Value = 667
$ws.Range("T174").Value = 15
